$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 19; rows 19-41 shift down to 20-42,
# keeping all of their existing values intact.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly record.
$ws.Cells.Item(19, 1).Value = 7
$ws.Cells.Item(19, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(19, 3).Value = "Ñuble"
$ws.Cells.Item(19, 4).Value = 44792
$ws.Cells.Item(19, 5).Value = 16
$ws.Cells.Item(19, 6).Value = 100112026
$ws.Cells.Item(19, 7).Value = "Haba"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 50
$ws.Cells.Item(19, 11).Value = 12000
$ws.Cells.Item(19, 12).Value = 12000
$ws.Cells.Item(19, 13).Value = 12000
$ws.Cells.Item(19, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(19, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(19, 16).Value = 480
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = "Hortaliza"
